$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Logistic Regression row (row 2) with new metric values
$ws.Range("B2").Value = 0.9033057851239669
$ws.Range("C2").Value = 0.9044857886687019
$ws.Range("D2").Value = 0.9033057851239669
$ws.Range("E2").Value = 0.9022504752972083

# Remove rows for Lasso, Support Vector Classifier, CART, Random Forest (rows 3-6)
$ws.Range("A3:E6").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# After the delete, LightGBM is now row 3 and XGBoost is now row 4.
# Remove the XGBoost row (now row 4)
$ws.Range("A4:E4").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

$wb.Save()
